# Update the cryptocurrency price/volume snapshot (columns D and E)
# to the latest scraped values, as produced by the scheduled
# GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.132.70"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "'2.367.62"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'0.698"
$ws.Range("E5").Value = "  +6.07%  "
$ws.Range("D6").Value = "'241.81"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").Value = "'77.00"
$ws.Range("E7").Value = "  +5.28%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +15.93%  "
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").Value = "'57.48"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "'33.47"
$ws.Range("E12").Value = "  +19.85%  "
$ws.Range("D13").Value = "'7.51"
$ws.Range("E13").Value = "  +13.35%  "
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "'2.721.66"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "'16.79"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("E17").Value = "  +5.64%  "
$ws.Range("D18").Value = "'2.364.29"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'44.064.44"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "'6.72"
$ws.Range("E21").Value = "  +6.48%  "
$ws.Range("D22").Value = "'77.90"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").Value = "'260.14"
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'3.75"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "'2.54"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("E27").Value = "  +17.08%  "
$ws.Range("D28").Value = "'10.93"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("E29").Value = "  +3.97%  "
$ws.Range("D30").Value = "'2.24"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'175.27"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("D34").Value = "'5.40"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").Value = "'0.0766"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("D36").Value = "'5.45"
$ws.Range("E36").Value = "  +7.25%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("D41").Value = "'0.222"
$ws.Range("E41").Value = "  +24.06%  "
$ws.Range("D42").Value = "'19.42"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "'9.21"
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("E44").Value = "  +12.82%  "
$ws.Range("D45").Value = "'4.90"
$ws.Range("E45").Value = "  +10.96%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'2.55"
$ws.Range("E47").Value = "  +11.42%  "
$ws.Range("D48").Value = "'1.26"
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D50").Value = "'102.58"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").Value = "'56.15"
$ws.Range("E51").Value = "  +9.16%  "

Write-Output "Updated 41 crypto rows (columns D/E) in-place."
